$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "About" sheet: add new notes under the "Notes:" section.
# ---------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# Insert 5 new rows right after the existing "Notes:" row (row 17),
# pushing the old note rows (18-20) down to 23-25.
$wsAbout.Rows("18:22").Insert()
# The inserted rows pick up the (bold) formatting of the row above;
# these are plain body-text notes, so drop back to the regular style.
$wsAbout.Range("A18:A22").ClearFormats()

$wsAbout.Range("A18").Value = "For vehicle types that can use electricity, this variable specifies the percentage"
$wsAbout.Range("A19").Value = "reduction in fuel use (on a BTU basis) relative to the typical fuel type for that vehicle"
$wsAbout.Range("A20").Value = "type (e.g. gasoline for LDVs, diesel for HDVs, etc.) due to the fact that electricity"
$wsAbout.Range("A21").Value = "can be converted into work more efficiently than other fuel types."
# Row 22 is left blank as a spacer (matches the source row layout).

# ---------------------------------------------------------------
# "PTFURfE" sheet: relabel the header, recompute aircraft/ships,
# and make the fuel-use-reduction column wider / taller to match.
# ---------------------------------------------------------------
$wsPT = $wb.Worksheets.Item("PTFURfE")

$wsPT.Range("A1").Value = "Percentage Reduction (dimensionless)"
$wsPT.Range("A1").WrapText = $true
$wsPT.Rows.Item(1).RowHeight = 45
$wsPT.Columns.Item(1).ColumnWidth = 15.95

# aircraft (row 4) and ships (row 6) now mirror rail (row 5) instead of
# being hard-coded zeros.
$wsPT.Range("B4").Formula = "=B5"
$wsPT.Range("C4").Formula = "=C5"
$wsPT.Range("B4").NumberFormat = "0.000"
$wsPT.Range("C4").NumberFormat = "0.000"

$wsPT.Range("B6").Formula = "=B5"
$wsPT.Range("C6").Formula = "=C5"
$wsPT.Range("B6").NumberFormat = "0.000"
$wsPT.Range("C6").NumberFormat = "0.000"

# ---------------------------------------------------------------
# Back on "About": append the aircraft/ships explanatory note.
# ---------------------------------------------------------------
$wsAbout.Range("A27").Value = "Aircraft and ships are assumed to be the same as rail, since they all use large engines"
$wsAbout.Range("A28").Value = "intended to move heavy craft."

Write-Host "edit complete"
